$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PC1 (column B) and PC2 (column C) values for rows 2-17
$data = @(
    @{ Row = 2;  B = -0.2857829782529795;  C = 0.1251746740930746 },
    @{ Row = 3;  B = -0.1713641216280674;  C = 0.2691628866934457 },
    @{ Row = 4;  B = -0.2672091298857421;  C = 0.08885964872219815 },
    @{ Row = 5;  B = -0.2037609987188569;  C = -0.3370909558608488 },
    @{ Row = 6;  B = -0.343731815768369;   C = -0.3524625223551583 },
    @{ Row = 7;  B = -0.3402266287114395;  C = -0.353055010291595 },
    @{ Row = 8;  B = -0.3793864944881054;  C = -0.009856245843924048 },
    @{ Row = 9;  B = -0.4017671282786786;  C = -0.1027945017834699 },
    @{ Row = 10; B = -0.3336618454222509;  C = 0.02231914608510188 },
    @{ Row = 11; B = -0.2269780937474418;  C = 0.4251283202836275 },
    @{ Row = 12; B = -0.06721544753267006; C = 0.2287364523425529 },
    @{ Row = 13; B = -0.187401532994686;   C = 0.4674712650727001 },
    @{ Row = 14; B = -0.1383694910781134;  C = 0.278486136237694 },
    @{ Row = 15; B = -0.04032849852316302; C = 0.006838485553286268 },
    @{ Row = 16; B = -0.07730361676731469; C = 0.01122341979271957 },
    @{ Row = 17; B = 0.08690005603000782;  C = -0.03171004826778088 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
